{"js": "// \"typo and wording fixes\"\n//\n// 1) \"...world shaper has grant them a cumulative...\" -> \"...has grants them a cumulative...\"\n// 2) Move the stray \"_GoBack\" bookmark from the \"Dichotomic Existence\" paragraph\n//    to the end of the \"...5% chance to ignore critical or precision damage\" paragraph.\n// 3) \"...shift around to keep you alive better than most, the negative hp...\" ->\n//    \"...shift around to compensate for damage, the negative hp...\"\n// 4) \"...treated as being affected by the grease spell\" ->\n//    \"...treated as being affected by the grease spell and leave a trail of grease\n//    where you go and on things you touch\"\n\nconst body = context.document.body;\n\n// --- 1) grant -> grants -------------------------------------------------\nconst grantScope = body.search(\"has grant them a cumulative\", { matchCase: true });\ngrantScope.load(\"items\");\nawait context.sync();\nconst grantWord = grantScope.items[0].search(\"grant\", { matchCase: true });\ngrantWord.load(\"items\");\nawait context.sync();\ngrantWord.items[0].getRange(\"End\").insertText(\"s\", Word.InsertLocation.replace);\n\n// --- 2) move the _GoBack bookmark ---------------------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nconst anomalyScope = body.search(\"5% chance to ignore critical or precision damage\", { matchCase: true });\nanomalyScope.load(\"items\");\nawait context.sync();\nanomalyScope.items[0].getRange(\"End\").insertBookmark(\"_GoBack\");\n\n// --- 3) reword the Redundancy ability description ------------------------\nconst redundancy = body.search(\n  \"Your organs and internals may shift around to keep you alive better than most, the negative hp required to kill you is double normal\",\n  { matchCase: true }\n);\nredundancy.load(\"items\");\nawait context.sync();\nredundancy.items[0].insertText(\n  \"Your organs and internals may shift around to compensate for damage, the negative hp required to kill you is double normal\",\n  Word.InsertLocation.replace\n);\n\n// --- 4) extend the Lubricating Fluid description --------------------------\nconst grease = body.search(\"treated as being affected by the grease spell\", { matchCase: true });\ngrease.load(\"items\");\nawait context.sync();\ngrease.items[0].getRange(\"End\").insertText(\n  \" and leave a trail of grease where you go and on things you touch\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# \"typo and wording fixes\"\n#\n# 1) \"...world shaper has grant them a cumulative...\" -> \"...has grants them a cumulative...\"\n# 2) Move the stray \"_GoBack\" bookmark from the \"Dichotomic Existence\" paragraph\n#    to the end of the \"...5% chance to ignore critical or precision damage\" paragraph.\n# 3) \"...shift around to keep you alive better than most, the negative hp...\" ->\n#    \"...shift around to compensate for damage, the negative hp...\"\n# 4) \"...treated as being affected by the grease spell\" ->\n#    \"...treated as being affected by the grease spell and leave a trail of grease\n#    where you go and on things you touch\"\n\n$d = $word.ActiveDocument\n\n# --- 1) grant -> grants --------------------------------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"has grant them a cumulative\"\n$find.Replacement.Text = \"has grants them a cumulative\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- 2) move the _GoBack bookmark ----------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate the end of \"...5% chance to ignore critical or precision damage\" and drop\n# a temporary marker there (Bookmarks.Add on a truly zero-length Range is unreliable\n# in this area of the document, so we bookmark a tiny run of real text instead and\n# then delete that text, which leaves the bookmark collapsed in the right spot).\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"5% chance to ignore critical or precision damage\"\n$found = $find.Execute()\nif ($found) {\n    $rng.Collapse(0)\n    $rng.InsertAfter(\"TempGoBackAnchor\")\n\n    $markerRng = $d.Content\n    $markerFind = $markerRng.Find\n    $markerFind.Text = \"TempGoBackAnchor\"\n    $markerFind.Execute() | Out-Null\n\n    $d.Bookmarks.Add(\"_GoBack\", $markerRng)\n    $markerRng.Text = \"\"\n}\n\n# --- 3) reword the Redundancy ability description -------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"Your organs and internals may shift around to keep you alive better than most, the negative hp required to kill you is double normal\"\n$find.Replacement.Text = \"Your organs and internals may shift around to compensate for damage, the negative hp required to kill you is double normal\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- 4) extend the Lubricating Fluid description ---------------------------\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"treated as being affected by the grease spell\"\n$find.Replacement.Text = \"treated as being affected by the grease spell and leave a trail of grease where you go and on things you touch\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n"}
